# Punchlist update — "Punched a few items out"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window / view cosmetics -------------------------------------------------
# Tab-bar split ratio (bookViews/workbookView@tabRatio 169 -> 84)
$win = $excel.ActiveWindow
$win.TabRatio = 84

# --- Remove the stale "Updated 12/17/14" note under the title ---------------
$ws.Range("A2").ClearContents()

# --- Mark the first few coding items as done, with dates --------------------
$ws.Range("E5").Value = "Done"
$ws.Range("F5").Value = 41995
$ws.Range("F5").NumberFormat = "MM/DD/YY"

$ws.Range("E6").Value = "Done – created Rpath object in ecopath()"
$ws.Range("F6").Value = 41995
$ws.Range("F6").NumberFormat = "MM/DD/YY"

$ws.Range("E7").Value = "First pass done"
$ws.Range("F7").Value = 42005
$ws.Range("F7").NumberFormat = "MM/DD/YY"

$ws.Range("F8").Value = 42005
$ws.Range("F8").NumberFormat = "MM/DD/YY"

$ws.Range("F10").Value = 42019
$ws.Range("F10").NumberFormat = "MM/DD/YY"

# --- Trim the trailing space on the RCPP description -------------------------
$ws.Range("C9").Value = "The R package RCPP integrates R and C++ directly, may gain efficiency in ecosim()"

# --- Add a note about hiring a programmer ------------------------------------
$ws.Range("E16").Value = "C code already exists? - Kerim may hire a programmer to do this next fall"

# --- Leave the cursor parked on A2, matching the saved selection ------------
$ws.Range("A2").Select()
